$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row, shifting existing rows 24-81 down to 25-82.
# Restrict the insert to the used columns (A:M) so we don't stamp
# formatting across all 16384 columns like a full Rows.Insert() would.
$ws.Range("A24:M24").Insert(-4121)

# Populate the newly inserted row 24 with the new indicator data (4.3.b)
$ws.Cells.Item(24, 1).Value = "Z04_B03_P01_Ib02_I01"
$ws.Cells.Item(24, 2).Value = "Z04_B03_P01_Ib02"
$ws.Cells.Item(24, 3).Value = "4.3.b"
$ws.Cells.Item(24, 4).Value = "Anteil der 17 bis 18-Jährigen mit (angestrebter) Studienberechtigung"
$ws.Cells.Item(24, 5).Value = "XXXAnteil der 17 bis 18-Jährigen mit (angestrebter) Studienberechtigung"
$ws.Cells.Item(24, 6).Value = "Anteil der 17 bis 18-Jährigen mit (angestrebter) Studienberechtigung"
$ws.Cells.Item(24, 7).Value = "XXXAnteil der 17 bis 18-Jährigen mit (angestrebter) Studienberechtigung"
$ws.Cells.Item(24, 8).Value = "Differenz zwischen den Gruppen „Keine Risikolage“ und „Mindestens eine Risikolage“ verringern bzw. angleichen"
$ws.Cells.Item(24, 9).Value = "XXXDifferenz zwischen den Gruppen „Keine Risikolage“ und „Mindestens eine Risikolage“ verringern bzw. angleichen"
$ws.Cells.Item(24, 10).Value = "Differenz zwischen den Gruppen „Keine Risikolage“ und „Mindestens eine Risikolage“ verringern bzw. angleichen"
$ws.Cells.Item(24, 11).Value = "XXXDifferenz zwischen den Gruppen „Keine Risikolage“ und „Mindestens eine Risikolage“ verringern bzw. angleichen"
$ws.Cells.Item(24, 12).Value = "Anteil der 17 bis 18-Jährigen mit (angestrebter) Studienberechtigung"
$ws.Cells.Item(24, 13).Value = "XXXAnteil der 17 bis 18-Jährigen mit (angestrebter) Studienberechtigung"

# Copy the style of the row above (row 23, a normal data row) into the new row
$ws.Range("A23:M23").Copy()
$ws.Range("A24:M24").PasteSpecial(-4122)
